# Apply updated LR-pair data per Dr Hou advice: recompute clusters with corrected
# sending/target cluster assignments and updated expression statistics,
# and extend the table from 10 to 15 data rows (M2 sending-cluster rows added).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Lrpap1'
$ws.Range("C2").Value = 'Vldlr'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.927655000000001
$ws.Range("H2").Value = 23.782965
$ws.Range("I2").Value = 0.1200556835465435
$ws.Range("J2").Value = 0.1200556835465435
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3081963333333333
$ws.Range("N2").Value = 0.924589
$ws.Range("O2").Value = 0.09210955608663024
$ws.Range("P2").Value = 0.09210955608663024
$ws.Range("Q2").Value = 2.443274202931667
$ws.Range("R2").Value = 21.989467826385
$ws.Range("S2").Value = 0.01105827571714908
$ws.Range("T2").Value = 0.01105827571714908

# Row 3
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Lrpap1'
$ws.Range("C3").Value = 'Vldlr'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.927655000000001
$ws.Range("H3").Value = 23.782965
$ws.Range("I3").Value = 0.1200556835465435
$ws.Range("J3").Value = 0.1200556835465435
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.207039333333333
$ws.Range("N3").Value = 6.621118
$ws.Range("O3").Value = 0.65961009678592
$ws.Range("P3").Value = 0.6596100967859201
$ws.Range("Q3").Value = 17.49664640609667
$ws.Range("R3").Value = 157.46981765487
$ws.Range("S3").Value = 0.07918994104383534
$ws.Range("T3").Value = 0.07918994104383534

# Row 4
$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Lrpap1'
$ws.Range("C4").Value = 'Vldlr'
$ws.Range("D4").Value = 'sCs'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.927655000000001
$ws.Range("H4").Value = 23.782965
$ws.Range("I4").Value = 0.1200556835465435
$ws.Range("J4").Value = 0.1200556835465435
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.83074
$ws.Range("N4").Value = 2.49222
$ws.Range("O4").Value = 0.2482803471274497
$ws.Range("P4").Value = 0.2482803471274497
$ws.Range("Q4").Value = 6.585820114700001
$ws.Range("R4").Value = 59.2723810323
$ws.Range("S4").Value = 0.02980746678555907
$ws.Range("T4").Value = 0.02980746678555907

# Row 5
$ws.Range("A5").Value = 'FAPs'
$ws.Range("B5").Value = 'Lrpap1'
$ws.Range("C5").Value = 'Vldlr'
$ws.Range("D5").Value = 'ECs'
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.980532
$ws.Range("H5").Value = 38.941596
$ws.Range("I5").Value = 0.1965759915205419
$ws.Range("J5").Value = 0.1965759915205419
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3081963333333333
$ws.Range("N5").Value = 0.924589
$ws.Range("O5").Value = 0.09210955608663024
$ws.Range("P5").Value = 0.09210955608663024
$ws.Range("Q5").Value = 4.000552367116001
$ws.Range("R5").Value = 36.00497130404401
$ws.Range("S5").Value = 0.0181065273162463
$ws.Range("T5").Value = 0.0181065273162463

# Row 6
$ws.Range("A6").Value = 'FAPs'
$ws.Range("B6").Value = 'Lrpap1'
$ws.Range("C6").Value = 'Vldlr'
$ws.Range("D6").Value = 'FAPs'
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.980532
$ws.Range("H6").Value = 38.941596
$ws.Range("I6").Value = 0.1965759915205419
$ws.Range("J6").Value = 0.1965759915205419
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.207039333333333
$ws.Range("N6").Value = 6.621118
$ws.Range("O6").Value = 0.65961009678592
$ws.Range("P6").Value = 0.6596100967859201
$ws.Range("Q6").Value = 28.648544691592
$ws.Range("R6").Value = 257.836902224328
$ws.Range("S6").Value = 0.1296635087926528
$ws.Range("T6").Value = 0.1296635087926528

# Row 7
$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Lrpap1'
$ws.Range("C7").Value = 'Vldlr'
$ws.Range("D7").Value = 'sCs'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.980532
$ws.Range("H7").Value = 38.941596
$ws.Range("I7").Value = 0.1965759915205419
$ws.Range("J7").Value = 0.1965759915205419
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.83074
$ws.Range("N7").Value = 2.49222
$ws.Range("O7").Value = 0.2482803471274497
$ws.Range("P7").Value = 0.2482803471274497
$ws.Range("Q7").Value = 10.78344715368
$ws.Range("R7").Value = 97.05102438312001
$ws.Range("S7").Value = 0.04880595541164275
$ws.Range("T7").Value = 0.04880595541164275

# Row 8
$ws.Range("A8").Value = 'M1'
$ws.Range("B8").Value = 'Lrpap1'
$ws.Range("C8").Value = 'Vldlr'
$ws.Range("D8").Value = 'ECs'
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 21.95152166666667
$ws.Range("H8").Value = 65.85456500000001
$ws.Range("I8").Value = 0.3324318400054526
$ws.Range("J8").Value = 0.3324318400054526
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.3081963333333333
$ws.Range("N8").Value = 0.924589
$ws.Range("O8").Value = 0.09210955608663024
$ws.Range("P8").Value = 0.09210955608663024
$ws.Range("Q8").Value = 6.76537848875389
$ws.Range("R8").Value = 60.88840639878501
$ws.Range("S8").Value = 0.03062014921196393
$ws.Range("T8").Value = 0.03062014921196393

# Row 9
$ws.Range("A9").Value = 'M1'
$ws.Range("B9").Value = 'Lrpap1'
$ws.Range("C9").Value = 'Vldlr'
$ws.Range("D9").Value = 'FAPs'
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 21.95152166666667
$ws.Range("H9").Value = 65.85456500000001
$ws.Range("I9").Value = 0.3324318400054526
$ws.Range("J9").Value = 0.3324318400054526
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.207039333333333
$ws.Range("N9").Value = 6.621118
$ws.Range("O9").Value = 0.65961009678592
$ws.Range("P9").Value = 0.6596100967859201
$ws.Range("Q9").Value = 48.44787174485223
$ws.Range("R9").Value = 436.0308457036701
$ws.Range("S9").Value = 0.2192753981607181
$ws.Range("T9").Value = 0.2192753981607181

# Row 10
$ws.Range("A10").Value = 'M1'
$ws.Range("B10").Value = 'Lrpap1'
$ws.Range("C10").Value = 'Vldlr'
$ws.Range("D10").Value = 'sCs'
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 21.95152166666667
$ws.Range("H10").Value = 65.85456500000001
$ws.Range("I10").Value = 0.3324318400054526
$ws.Range("J10").Value = 0.3324318400054526
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.83074
$ws.Range("N10").Value = 2.49222
$ws.Range("O10").Value = 0.2482803471274497
$ws.Range("P10").Value = 0.2482803471274497
$ws.Range("Q10").Value = 18.23600710936667
$ws.Range("R10").Value = 164.1240639843
$ws.Range("S10").Value = 0.0825362926327706
$ws.Range("T10").Value = 0.08253629263277061

# Row 11
$ws.Range("A11").Value = 'M2'
$ws.Range("B11").Value = 'Lrpap1'
$ws.Range("C11").Value = 'Vldlr'
$ws.Range("D11").Value = 'ECs'
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 19.402266
$ws.Range("H11").Value = 58.20679800000001
$ws.Range("I11").Value = 0.293826144929599
$ws.Range("J11").Value = 0.293826144929599
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.3081963333333333
$ws.Range("N11").Value = 0.924589
$ws.Range("O11").Value = 0.09210955608663024
$ws.Range("P11").Value = 0.09210955608663024
$ws.Range("Q11").Value = 5.979707239558
$ws.Range("R11").Value = 53.81736515602201
$ws.Range("S11").Value = 0.02706419577611125
$ws.Range("T11").Value = 0.02706419577611125

# Row 12
$ws.Range("A12").Value = 'M2'
$ws.Range("B12").Value = 'Lrpap1'
$ws.Range("C12").Value = 'Vldlr'
$ws.Range("D12").Value = 'FAPs'
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 19.402266
$ws.Range("H12").Value = 58.20679800000001
$ws.Range("I12").Value = 0.293826144929599
$ws.Range("J12").Value = 0.293826144929599
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.207039333333333
$ws.Range("N12").Value = 6.621118
$ws.Range("O12").Value = 0.65961009678592
$ws.Range("P12").Value = 0.6596100967859201
$ws.Range("Q12").Value = 42.821564217796
$ws.Range("R12").Value = 385.3940779601641
$ws.Range("S12").Value = 0.1938106918952466
$ws.Range("T12").Value = 0.1938106918952466

# Row 13
$ws.Range("A13").Value = 'M2'
$ws.Range("B13").Value = 'Lrpap1'
$ws.Range("C13").Value = 'Vldlr'
$ws.Range("D13").Value = 'sCs'
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 19.402266
$ws.Range("H13").Value = 58.20679800000001
$ws.Range("I13").Value = 0.293826144929599
$ws.Range("J13").Value = 0.293826144929599
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.83074
$ws.Range("N13").Value = 2.49222
$ws.Range("O13").Value = 0.2482803471274497
$ws.Range("P13").Value = 0.2482803471274497
$ws.Range("Q13").Value = 16.11823845684
$ws.Range("R13").Value = 145.06414611156
$ws.Range("S13").Value = 0.07295125725824118
$ws.Range("T13").Value = 0.0729512572582412

# Row 14
$ws.Range("A14").Value = 'sCs'
$ws.Range("B14").Value = 'Lrpap1'
$ws.Range("C14").Value = 'Vldlr'
$ws.Range("D14").Value = 'ECs'
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3.771175666666667
$ws.Range("H14").Value = 11.313527
$ws.Range("I14").Value = 0.05711033999786299
$ws.Range("J14").Value = 0.05711033999786299
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.3081963333333333
$ws.Range("N14").Value = 0.924589
$ws.Range("O14").Value = 0.09210955608663024
$ws.Range("P14").Value = 0.09210955608663024
$ws.Range("Q14").Value = 1.162262512822556
$ws.Range("R14").Value = 10.460362615403
$ws.Range("S14").Value = 0.005260408065159683
$ws.Range("T14").Value = 0.005260408065159683

# Row 15
$ws.Range("A15").Value = 'sCs'
$ws.Range("B15").Value = 'Lrpap1'
$ws.Range("C15").Value = 'Vldlr'
$ws.Range("D15").Value = 'FAPs'
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3.771175666666667
$ws.Range("H15").Value = 11.313527
$ws.Range("I15").Value = 0.05711033999786299
$ws.Range("J15").Value = 0.05711033999786299
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.207039333333333
$ws.Range("N15").Value = 6.621118
$ws.Range("O15").Value = 0.65961009678592
$ws.Range("P15").Value = 0.6596100967859201
$ws.Range("Q15").Value = 8.323133029242889
$ws.Range("R15").Value = 74.90819726318601
$ws.Range("S15").Value = 0.0376705568934672
$ws.Range("T15").Value = 0.03767055689346721

# Row 16
$ws.Range("A16").Value = 'sCs'
$ws.Range("B16").Value = 'Lrpap1'
$ws.Range("C16").Value = 'Vldlr'
$ws.Range("D16").Value = 'sCs'
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 3.771175666666667
$ws.Range("H16").Value = 11.313527
$ws.Range("I16").Value = 0.05711033999786299
$ws.Range("J16").Value = 0.05711033999786299
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.83074
$ws.Range("N16").Value = 2.49222
$ws.Range("O16").Value = 0.2482803471274497
$ws.Range("P16").Value = 0.2482803471274497
$ws.Range("Q16").Value = 3.132866473326667
$ws.Range("R16").Value = 28.19579825994
$ws.Range("S16").Value = 0.0141793750392361
$ws.Range("T16").Value = 0.0141793750392361

